# RTxManualControlBOM.xlsx edit script
# Adds a new "Order Qty" column (C) with quantity-to-order numbers, fills in
# "Order PN" (Digikey part numbers) for every BOM line, and adds a parallel
# "Female-*" connector line after every pin-header line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column C ("Order Qty") - shifts old C..I to D..J
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()
$ws.Columns(3).ColumnWidth = 5.6

# ---------------------------------------------------------------------------
# 2. Insert the six new "Female-*" rows at their final row numbers.
#    Doing these top-to-bottom at the *final* row number works because each
#    insertion only affects rows below it, which have not been processed yet.
# ---------------------------------------------------------------------------
$ws.Rows(9).Insert()
$ws.Rows(9).RowHeight = 17.25

$ws.Rows(11).Insert()
$ws.Rows(11).RowHeight = 17.25

$ws.Rows(13).Insert()
$ws.Rows(13).RowHeight = 17.25

$ws.Rows(15).Insert()
$ws.Rows(15).RowHeight = 17.25

$ws.Rows(17).Insert()
$ws.Rows(17).RowHeight = 17.25

# Row 19 is brand new territory past the old data, no insert needed, and it
# keeps the default row height (no explicit ht in the target).

# Header row is now taller (wraps to two lines with the new column)
$ws.Rows(1).RowHeight = 29.25

# ---------------------------------------------------------------------------
# 3. Fill in header row
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,3).Value = "Order Qty"

# ---------------------------------------------------------------------------
# 4. Fill in "Order Qty" (col C) and "Order PN" (col D) for every line, plus
#    the new Female-* rows in col D/F.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,3).Value = 21
$ws.Cells.Item(2,4).Value = "311-1376-1-ND"

$ws.Cells.Item(3,3).Value = 10
$ws.Cells.Item(3,4).Value = "311-1181-1-ND"

$ws.Cells.Item(4,3).Value = 10
$ws.Cells.Item(4,4).Value = "311-10.0KFRCT-ND"

$ws.Cells.Item(5,3).Value = 10
$ws.Cells.Item(5,4).Value = "311-1.00KFRCT-ND"

$ws.Cells.Item(6,3).Value = 18
$ws.Cells.Item(6,4).Value = "311-280FRCT-ND"

$ws.Cells.Item(7,3).Value = 12
$ws.Cells.Item(7,4).Value = "311-620FRCT-ND"

$ws.Cells.Item(8,3).Value = 6
$ws.Cells.Item(8,4).Value = "609-3406-ND"

$ws.Cells.Item(9,4).Value = "A106652-ND"
$ws.Cells.Item(9,6).Value = "Female-1X6"

$ws.Cells.Item(10,3).Value = 6
$ws.Cells.Item(10,4).Value = "609-3406-ND"

$ws.Cells.Item(11,3).Value = 6
$ws.Cells.Item(11,6).Value = "Female-1X8"

$ws.Cells.Item(12,3).Value = 3

$ws.Cells.Item(13,3).Value = 3
$ws.Cells.Item(13,6).Value = "Female-1X2"

$ws.Cells.Item(14,3).Value = 6

$ws.Cells.Item(15,3).Value = 6
$ws.Cells.Item(15,6).Value = "Female-1X6 Polarized"

$ws.Cells.Item(16,3).Value = 6

$ws.Cells.Item(17,3).Value = 6
$ws.Cells.Item(17,6).Value = "Female-2X3 Polarized"

$ws.Cells.Item(18,3).Value = 3

$ws.Cells.Item(19,3).Value = 3
$ws.Cells.Item(19,6).Value = "Female-1X3 Polarized"

# Row 19 cells came from plain new-row territory (no Insert was used) so they
# do not carry the "wrap text / Times New Roman" look of the rest of the
# table yet. Copy that formatting over (only the two cells actually used)
# from row 18 without touching values.
$ws.Range("C18").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("F18").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Selection shown in the saved file
# ---------------------------------------------------------------------------
$ws.Range("I10").Select() | Out-Null
